$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "305.01"
Set-TextValue "E2" "0.26%"
Set-TextValue "D3" "37.07"
Set-TextValue "E3" "6.20%"
Set-TextValue "D4" "5.011"
Set-TextValue "E4" "-3.06%"
Set-TextValue "D5" "0.07900"
Set-TextValue "E5" "0.98%"
Set-TextValue "D6" "2.208"
Set-TextValue "E6" "-4.69%"
Set-TextValue "D7" "8.026"
Set-TextValue "E7" "-0.51%"
Set-TextValue "D9" "0.9211"
Set-TextValue "E9" "-0.49%"
Set-TextValue "D10" "0.09671"
Set-TextValue "E10" "-4.12%"
Set-TextValue "D11" "0.1892"
Set-TextValue "E11" "3.46%"
Set-TextValue "D12" "0.08622"
Set-TextValue "E12" "0.71%"
Set-TextValue "D13" "0.03684"
Set-TextValue "E13" "8.87%"
Set-TextValue "D14" "0.09990"
Set-TextValue "E14" "0.87%"
Set-TextValue "E15" "-1.23%"
Set-TextValue "D16" "0.005631"
Set-TextValue "E16" "-1.83%"
Set-TextValue "D17" "3.476"
Set-TextValue "E17" "0.05%"
Set-TextValue "E18" "6.96%"
Set-TextValue "E19" "-0.09%"
Set-TextValue "D20" "0.1317"
Set-TextValue "E20" "-0.72%"
Set-TextValue "E21" "4.63%"
Set-TextValue "D22" "0.2199"
Set-TextValue "E22" "-3.17%"
Set-TextValue "D23" "0.04557"
Set-TextValue "E23" "-2.01%"
Set-TextValue "D24" "0.001234"
Set-TextValue "E24" "1.37%"
Set-TextValue "E25" "3.13%"
Set-TextValue "D26" "0.0001401"
Set-TextValue "E26" "7.70%"
Set-TextValue "E27" "39.78%"
Set-TextValue "D39" "0.01843"
Set-TextValue "E39" "5.01%"
Set-TextValue "D40" "0.04758"
Set-TextValue "E40" "0.24%"
Set-TextValue "D41" "0.008135"
Set-TextValue "E41" "4.60%"
Set-TextValue "D42" "0.1399"
Set-TextValue "E42" "-0.86%"
Set-TextValue "D44" "0.002231"
Set-TextValue "E44" "-2.62%"
Set-TextValue "D45" "0.01006"
Set-TextValue "E45" "0.95%"
Set-TextValue "D46" "0.00006267"
Set-TextValue "E46" "3.48%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.08%"
Set-TextValue "D48" "0.0005798"
Set-TextValue "E48" "-0.05%"
Set-TextValue "D49" "37.49"
Set-TextValue "E49" "866.25%"
Set-TextValue "E50" "-36.06%"
Set-TextValue "D51" "0.00002099"
Set-TextValue "E51" "-0.08%"
